# Update the "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages update).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row => New F-column value, for sheet "展览"
$updates1 = @{
    6  = 1094
    7  = 1434
    11 = 71
    12 = 164
    15 = 1352
    16 = 113
    17 = 103
    18 = 276
    20 = 29
    21 = 654
    24 = 221
    26 = 5855
    27 = 64
    29 = 98
    31 = 14498
    32 = 1437
    36 = 8561
    37 = 621
    38 = 4207
    39 = 139
}

foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Row => New F-column value, for sheet "全部类型"
$updates4 = @{
    6  = 1094
    7  = 1434
    11 = 71
    12 = 164
    15 = 1352
    16 = 113
    17 = 103
    18 = 276
    21 = 29
    22 = 654
    26 = 221
    29 = 5855
    30 = 64
    32 = 98
    34 = 14498
    35 = 1437
    39 = 8561
    40 = 621
    41 = 4207
    42 = 139
}

foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
